# Apply the "Artfynd" sheet update: rows 19-30 get their species-
# observation records re-synced (ids, counts, units, coordinates, etc.)
# to match the latest export from the source system.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = 111926622
$ws.Range("B19").Value = 90806
$ws.Range("I19").ClearContents()
$ws.Range("J19").ClearContents()
$ws.Range("Q19").Value = 663452
$ws.Range("R19").Value = 6602676

# Row 20
$ws.Range("A20").Value = 111926769
$ws.Range("B20").Value = 90806
$ws.Range("Q20").Value = 663476
$ws.Range("R20").Value = 6602651

# Row 21
$ws.Range("B21").Value = 90803

# Row 22
$ws.Range("A22").Value = 111927215
$ws.Range("B22").Value = 90806
$ws.Range("I22").Value = "10"
$ws.Range("J22").Value = "fruktkroppar"
$ws.Range("Q22").Value = 663486
$ws.Range("R22").Value = 6602647

# Row 23
$ws.Range("A23").Value = 112083905
$ws.Range("B23").Value = 98980
$ws.Range("I23").Value = "400"
$ws.Range("J23").Value = "stjälkar/strån/skott"
$ws.Range("Q23").Value = 663568
$ws.Range("R23").Value = 6602721

# Row 24
$ws.Range("A24").Value = 112084535
$ws.Range("AB24").Value = "11:46"
$ws.Range("AC24").Value = "Sötaktig mild smak (ej bitter)"
$ws.Range("B24").Value = 89047
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 3286
$ws.Range("F24").Value = "Flattoppad klubbsvamp"
$ws.Range("G24").Value = "Clavariadelphus truncatus"
$ws.Range("H24").Value = "(Quél.) Donk"
$ws.Range("I24").Value = "80"
$ws.Range("J24").Value = "fruktkroppar"
$ws.Range("K24").ClearContents()
$ws.Range("Q24").Value = 663374
$ws.Range("R24").Value = 6602611
$ws.Range("S24").Value = 10
$ws.Range("Z24").Value = "11:46"

# Row 25
$ws.Range("A25").Value = 112084040
$ws.Range("B25").Value = 98980
$ws.Range("I25").ClearContents()
$ws.Range("J25").ClearContents()
$ws.Range("K25").ClearContents()
$ws.Range("Q25").Value = 663585
$ws.Range("R25").Value = 6602704
$ws.Range("S25").Value = 10

# Row 26
$ws.Range("A26").Value = 112084114
$ws.Range("AB26").ClearContents()
$ws.Range("AC26").ClearContents()
$ws.Range("B26").Value = 98980
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 222498
$ws.Range("F26").Value = "Blåsippa"
$ws.Range("G26").Value = "Hepatica nobilis"
$ws.Range("H26").Value = "Schreb."
$ws.Range("I26").Value = "20"
$ws.Range("J26").Value = "plantor/tuvor"
$ws.Range("K26").Value = "fullt utvecklade blad"
$ws.Range("Q26").Value = 663577
$ws.Range("R26").Value = 6602715
$ws.Range("S26").Value = 5
$ws.Range("Z26").ClearContents()

# Row 27
$ws.Range("B27").Value = 98980

# Row 28
$ws.Range("A28").Value = 112083958
$ws.Range("B28").Value = 98980
$ws.Range("I28").Value = "10"
$ws.Range("J28").Value = "plantor/tuvor"
$ws.Range("K28").Value = "fullt utvecklade blad"
$ws.Range("Q28").Value = 663551
$ws.Range("R28").Value = 6602700
$ws.Range("S28").Value = 5

# Row 29
$ws.Range("A29").Value = 112083804
$ws.Range("B29").Value = 98980
$ws.Range("I29").Value = "300"
$ws.Range("J29").Value = "stjälkar/strån/skott"
$ws.Range("Q29").Value = 663572
$ws.Range("R29").Value = 6602738

# Row 30
$ws.Range("A30").Value = 112083991
$ws.Range("B30").Value = 98980
$ws.Range("Q30").Value = 663568
$ws.Range("R30").Value = 6602664
